# Append a freshly-scraped job listing (new row 11) and refresh every
# "fetched at" timestamp in column A to the new scrape time.
#
# Source workbook lists 22 jobs (rows 2-23) all stamped
# "2025-12-01 18:39:07". The new scrape adds one more job
# ("管理システムの開発", work/detail/5445265) which sorts in right after
# row 10 (=> becomes the new row 11, pushing the old rows 11-23 down to
# 12-24), and every row - old and new alike - gets re-stamped with the
# new scrape time "2025-12-02 01:21:59".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$lastRow = 23
$newRow = 11

# Make room for the new job at row 11; this shifts old rows 11..23 down
# to 12..24 and copies row 10's formatting (incl. the Hyperlink style on
# column F) down onto the newly-opened row.
$ws.Rows.Item($newRow).Insert()
$lastRow = $lastRow + 1

# Fill in the newly inserted row.
$ws.Range("A11").Value = "2025-12-02 01:21:59"
$ws.Range("B11").Value = "管理システムの開発"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5445265"
$ws.Range("G11").Value = 103
$ws.Range("H11").Value = "◆開発 ◇管理"

# Every row's "fetched at" stamp (column A) refreshes to the new scrape
# timestamp, old rows and the newly-inserted one alike.
$ws.Range("A2:A$lastRow").Value = "2025-12-02 01:21:59"

# Row-insert doesn't renumber the worksheet's <hyperlinks> table, so the
# existing entries now point at the wrong rows (they still reference the
# pre-insert row numbers / urls). Rebuild the hyperlink for every URL
# cell in column F from scratch, in row order, so rIds realign with the
# final row layout.
$ws.Hyperlinks.Delete()
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
